# Append scrape: 2025-09-09 01:15 JST
# Two new job postings were picked up by the scraper and inserted (by
# priority score) into the "ランサーズ" sheet, while every existing row's
# "取得日時" (retrieved-at) timestamp was refreshed to the new run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")
if ($ws -eq $null) {
    $ws = $wb.ActiveSheet
}

$newTimestamp = "2025-09-09 01:15:32"

# --- 1. Remove existing hyperlinks -----------------------------------
# Row inserts below do not renumber the worksheet's <hyperlinks> table in
# this environment, so drop them now and rebuild the full, correctly
# numbered set once all rows are in their final place.
$ws.Hyperlinks.Delete()

# --- 2. Make room for the two new rows --------------------------------
# New row 4: "AIクローン制作" (score 303, between existing 305 and 298)
$ws.Rows.Item(4).Insert()
# New row 12 (after the first insert shifted everything below row 4 down
# by one, the old row 11 is now row 12): "仮想通貨EA" (score 63, between
# existing 93 and 50)
$ws.Rows.Item(12).Insert()

# --- 3. Fill in the two newly inserted rows ---------------------------
$ws.Range("A4").Value = $newTimestamp
$ws.Range("B4").Value = "【AIクローン制作】肉声・映像応答のプロデュース依頼"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5389313"
$ws.Range("G4").Value = 303
$ws.Range("H4").Value = "🔥AI,Ai"

$ws.Range("A12").Value = $newTimestamp
$ws.Range("B12").Value = "【仮想通貨】自動売買EAの開発依頼"
$ws.Range("C12").Value = "システム開発"
$ws.Range("D12").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E12").Value = "期限情報なし"
$ws.Range("F12").Value = "https://www.lancers.jp/work/detail/5389714"
$ws.Range("G12").Value = 63
$ws.Range("H12").Value = "◆開発"

# --- 4. Refresh the retrieval timestamp on every data row -------------
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# --- 5. Rebuild hyperlinks for column F, rows 2-18 ---------------------
for ($r = 2; $r -le 18; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $url = $cell.Value()
    $ws.Hyperlinks.Add($cell, $url) | Out-Null
}

Write-Host "Done updating ランサーズ sheet"
